# Auto-generated edit script: updates cryptocurrency price/volume data
# and re-orders a handful of rows per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.316.83'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '2.971.96'
$ws.Range("E3").Value = '  +2.34%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '199.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.97'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.548'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("E9").Value = '  +3.10%  '
$ws.Range("D10").Value = '2.971.61'
$ws.Range("E10").Value = '  +2.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.433'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.39%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").Value = '3.517.32'
$ws.Range("E14").Value = '  +2.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.98'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +6.31%  '
$ws.Range("D16").Value = '76.177.42'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000186'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").Value = '2.962.52'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.37'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +6.35%  '
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '373.31'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.79%  '
$ws.Range("E22").Value = '  +3.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.24'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.79'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.41%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.123.37'
$ws.Range("E25").Value = '  +3.04%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +2.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.62'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000105'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.996'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.25'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +7.45%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '513.95'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.58%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.39'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.93'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.24'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.46'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.383'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +12.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.97'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.106'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +16.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.111'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '181.90'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.80'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.92%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.91'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.56%  '
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.21'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.695'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +6.76%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.583'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.54%  '
$ws.Range("E50").Value = '  -2.28%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.81'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.68%  '
